$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 207, shifting existing rows 207-214 down to 208-215.
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row 207 with the new record's data.
$ws.Range("A207").Value2 = 8
$ws.Range("B207").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C207").Value2 = "Coquimbo"
$ws.Range("D207").Value2 = 44753
$ws.Range("E207").Value2 = 4
$ws.Range("F207").Value2 = 100112037
$ws.Range("G207").Value2 = "Cebollín"
$ws.Range("H207").Value2 = "Sin especificar"
$ws.Range("I207").Value2 = "Primera"
$ws.Range("J207").Value2 = 2000
$ws.Range("K207").Value2 = 1400
$ws.Range("L207").Value2 = 1600
$ws.Range("M207").Value2 = 1500
$ws.Range("N207").Value2 = "`$/paquete 6 unidades"
$ws.Range("O207").Value2 = "Provincia del Elquí"
$ws.Range("P207").Value2 = 250
$ws.Range("Q207").Value2 = 6
$ws.Range("R207").Value2 = "Hortaliza"
